$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts Type/Priority/Communication_Routes/
# Profile_Picture/Last_Contacted/Last_Meeting/Industry/Comments one column
# to the right, E..L).
$ws.Columns("D:D").Insert()

# New column D: "Phone_Number" header + one text value per data row.
# Force text format so the leading zero is preserved instead of Excel
# re-interpreting the value as a number.
$ws.Range("D1:D19").NumberFormat = "@"

$ws.Range("D1").Value = "Phone_Number"

$phones = @(
    "0126906297",
    "0126906298",
    "0126906299",
    "0126906300",
    "0126906301",
    "0126906302",
    "0126906303",
    "0126906304",
    "0126906305",
    "0126906306",
    "0126906307",
    "0126906308",
    "0126906309",
    "0126906310",
    "0126906311",
    "0126906312",
    "0126906313",
    "0126906314"
)

for ($i = 0; $i -lt $phones.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $phones[$i]
}

# Widen the new column to fit its content. Excel stores column widths in
# the OOXML "width" (MDW units) but the COM ColumnWidth property is in
# characters; character-width 16 + 1/6 round-trips to width="17".
$ws.Columns("D:D").ColumnWidth = 16 + 1/6

# Match the author's final selection/cursor position.
$ws.Range("D10").Select() | Out-Null
